$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text before writing, since many values look numeric
# (e.g. "1.00", "5.28") and would otherwise be auto-converted by Excel, losing
# their original trailing-zero text formatting / exact representation.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '64.614.30'
$ws.Range('E2').Value = '  +2.88%  '
$ws.Range('D3').Value = '2.526.84'
$ws.Range('E3').Value = '  +2.61%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '579.61'
$ws.Range('E5').Value = '  +0.87%  '
$ws.Range('D6').Value = '152.52'
$ws.Range('E6').Value = '  +3.88%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '0.539'
$ws.Range('E8').Value = '  +0.76%  '
$ws.Range('D9').Value = '2.526.72'
$ws.Range('E9').Value = '  +2.59%  '
$ws.Range('E10').Value = '  +0.42%  '
$ws.Range('E11').Value = '  -1.39%  '
$ws.Range('D12').Value = '5.28'
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('D13').Value = '0.354'
$ws.Range('E13').Value = '  -0.79%  '
$ws.Range('D14').Value = '29.21'
$ws.Range('E14').Value = '  +0.76%  '
$ws.Range('D15').Value = '0.0000180'
$ws.Range('E15').Value = '  +1.11%  '
$ws.Range('D16').Value = '2.986.23'
$ws.Range('E16').Value = '  +2.65%  '
$ws.Range('D17').Value = '64.134.31'
$ws.Range('E17').Value = '  +2.27%  '
$ws.Range('D18').Value = '2.524.11'
$ws.Range('E18').Value = '  +2.25%  '
$ws.Range('D19').Value = '7.85'
$ws.Range('E19').Value = '  -1.16%  '
$ws.Range('D20').Value = '10.93'
$ws.Range('E20').Value = '  -0.62%  '
$ws.Range('D21').Value = '4.26'
$ws.Range('E21').Value = '  +3.15%  '
$ws.Range('D22').Value = '328.98'
$ws.Range('E22').Value = '  +0.88%  '
$ws.Range('D23').Value = '2.23'
$ws.Range('E23').Value = '  +0.59%  '
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('D25').Value = '10.10'
$ws.Range('E25').Value = '  +0.61%  '
$ws.Range('D26').Value = '65.74'
$ws.Range('E26').Value = '  +0.31%  '
$ws.Range('D27').Value = '646.76'
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('D28').Value = '0.0000104'
$ws.Range('E28').Value = '  +5.33%  '
$ws.Range('E29').Value = '  +2.72%  '
$ws.Range('E30').Value = '  +4.82%  '
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.37%  '
$ws.Range('D32').Value = '8.04'
$ws.Range('E32').Value = '  +0.75%  '
$ws.Range('E33').Value = '  +1.39%  '
$ws.Range('E34').Value = '  +1.87%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('E36').Value = '  +0.60%  '
$ws.Range('D37').Value = '4.82'
$ws.Range('E37').Value = '  +1.49%  '
$ws.Range('D38').Value = '5.52'
$ws.Range('E38').Value = '  +2.67%  '
$ws.Range('D39').Value = '154.58'
$ws.Range('E39').Value = '  +1.64%  '
$ws.Range('B40').Value = 'EthereumClassic'
$ws.Range('C40').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D40').Value = '18.95'
$ws.Range('E40').Value = '  +1.35%  '
$ws.Range('B41').Value = 'PolygonEcosystemToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D41').Value = '0.372'
$ws.Range('E41').Value = '  +0.77%  '
$ws.Range('D42').Value = '2.81'
$ws.Range('E42').Value = '  -0.39%  '
$ws.Range('D43').Value = '1.79'
$ws.Range('E43').Value = '  +3.16%  '
$ws.Range('D44').Value = '163.05'
$ws.Range('E44').Value = '  +7.21%  '
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').Value = '0.0₆0300'
$ws.Range('E46').Value = '  -2.77%  '
$ws.Range('D47').Value = '15.53'
$ws.Range('E47').Value = '  +1.74%  '
$ws.Range('D48').Value = '3.64'
$ws.Range('E48').Value = '  +1.45%  '
$ws.Range('D49').Value = '21.37'
$ws.Range('E49').Value = '  +4.01%  '
$ws.Range('D50').Value = '0.619'
$ws.Range('E50').Value = '  +1.96%  '
$ws.Range('E51').Value = '  +1.20%  '

# Restore default (Normal) style on column D now that the text values are set,
# so no stray explicit text-number-format style is left applied to the cells.
$ws.Range('D2:D51').Style = 'Normal'
